# Update the picture "Picture 4" on slide 6 ("Showing off the menu"):
#   - crop a sliver off the right edge of the image (a:srcRect r="2682")
#   - shrink/reposition the picture frame to match the new crop
#     (off x 1371600 -> 1371601 EMU, ext cx 4854575 -> 4724400 EMU)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)

# Crop ~13.98pt (2.682%) off the right edge of the source image.
$shape.PictureFormat.CropRight = 13.979925

# Resize/reposition the picture frame to the new (post-crop) footprint.
# Values are expressed in points (EMU / 12700) as PowerPoint's object model expects.
$shape.Left = 108.0001
$shape.Width = 372
